$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.190.89"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.285.61"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("E14").Value = "  +17.81%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "2.630.88"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "2.297.33"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "43.133.32"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "255.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0903"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.77%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  -7.08%  "
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "108.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.23%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.76%  "
